$wb = $excel.ActiveWorkbook

# --- Sheet1: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "Omar Al Awani"
$summary.Range("B4").Value = 7724.94
$summary.Range("B6").Value = 4438
$summary.Range("B7").Value = 42793
$summary.Range("B8").Value = -38355
$summary.Range("B9").Value = 0.1

# --- Sheet2: Assets ---
$assets = $wb.Worksheets.Item("Assets")
# Delete the two "Vehicles" rows (rows 2 and 3), shifting remaining rows up
$assets.Rows.Item(2).Delete()
$assets.Rows.Item(2).Delete()
# Update the remaining data row (now row 2) and total row (now row 3)
$assets.Range("C2").Value = 4438
$assets.Range("C3").Value = 4438

# --- Sheet3: Liabilities ---
$liabilities = $wb.Worksheets.Item("Liabilities")
# Delete the "Auto Loans" and "Personal Loans" rows (rows 2 and 3), shifting remaining rows up
$liabilities.Rows.Item(2).Delete()
$liabilities.Rows.Item(2).Delete()
# Update the remaining data row (now row 2) and total row (now row 3)
$liabilities.Range("C2").Value = 42793
$liabilities.Range("D2").Value = 2140
$liabilities.Range("E2").Value = 1
$liabilities.Range("C3").Value = 42793
